# Update db playercharacter atk data
# Player Attack Data sheet (4th sheet): rename spear_JumpX_Attack -> spear_Jump_X_Attack
# and add 5 new jump attack variants, then resort the table descending by attackType.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$lo = $ws.ListObjects.Item(1)

# Grow the table from 16 data rows (17 incl header) to 21 data rows (22 incl header)
# to make room for 5 new attack rows.
$lo.Resize($ws.Range("A1:H22"))

# Row 7 currently holds spear_JumpX_Attack -- rename it in place (values unchanged).
$ws.Range("A7").Value = "spear_Jump_X_Attack"

# Tweak attackMultiply for the two vertical spear attacks.
$ws.Range("F8").Value = 0.4    # spear_Y_Attack: 0.5 -> 0.4
$ws.Range("F9").Value = 1.5    # spear_YUp_Attack: 3 -> 1.5

# Fill the 5 newly added rows (18-22) with the new jump attack variants.
$ws.Range("A18").Value = "spear_Jump_XX_Attack"
$ws.Range("B18").Value = -0.25
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1.5
$ws.Range("F18").Value = 1.5
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 3

$ws.Range("A19").Value = "spear_Jump_XXX_Attack"
$ws.Range("B19").Value = -0.5
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 1.5
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 5

$ws.Range("A20").Value = "spear_Jump_Y_Attack"
$ws.Range("B20").Value = -0.25
$ws.Range("C20").Value = -0.15
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 2.5
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 5

$ws.Range("A21").Value = "spear_Jump_Down_X_Attack"
$ws.Range("B21").Value = 0.8
$ws.Range("C21").Value = -0.4
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1.5
$ws.Range("F21").Value = 1.5
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 3

$ws.Range("A22").Value = "spear_Jump_Up_X_Attack"
$ws.Range("B22").Value = 0.8
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1.5
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 2

# Re-sort the table descending by attackType (column A), matching the author's
# re-organization of the attack list.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("A1:A22"), 0, 2)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Cosmetic adjustments that accompanied the data edit.
$ws.Columns.Item(1).ColumnWidth = 27.5
$ws.Range("H14").Select()
